$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (new diary entry, 1/29/2020) gets filled in first (date/time/
#     participants/goal), matching the order the author actually typed
#     things in (this is what produces the exact shared-string ordering
#     seen in the saved workbook). ---
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 43859

$ws.Range("B16").Value = "9/00 - 12:00 p.m."
$ws.Range("C16").Value = "Zeyu Huang, Yue Zhang"
$ws.Range("D16").Value = "Export uml graph"

# --- Then the author went back and filled in row 15 (1/25/2020 entry)
#     completely. ---
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 43855

$ws.Range("B15").Value = "2:00 - 5:00 p.m."
$ws.Range("C15").Value = "Zeyu Huang, Yue Zhang"
$ws.Range("D15").Value = "Learn about the features of elasticsearch and decide which two features we want to dive in"
$ws.Range("E15").Value = "We choose ip filtering and rank evaluation "
$ws.Range("G15").Value = "Confused "
$ws.Range("F15").Value = "Elasticsearch is an incredible project with tons of great features. It seems pretty intimidating at first glance"

# --- Then back to row 16 to finish it off. ---
$ws.Range("F16").Value = "The project has very loose coupling, so we don't relly need the uml graph for the whole project, just the parts we need"
$ws.Range("E16").Value = "Elasticsearch is too big to draw an uml graph, with over 10k java classes, the graph is impossible to print out, so we took the parts we looked into and drawed some smaller uml graph"
$ws.Range("G16").Value = "Good"

# Row heights grew to fit the newly-wrapped text.
$ws.Rows.Item(15).RowHeight = 51
$ws.Rows.Item(16).RowHeight = 85

# Final selection left on I16, as in the saved file.
$ws.Range("I16").Select() | Out-Null
